{"js": "// Tiny update to readme:\n//   1. \"...alongside our public with the following authors:\"\n//        -> \"...alongside our manuscript with the following authors:\"\n//   2. \"Version 0.1 updated on 02/06/17\" -> \"Version 0.1 updated on 04/27/18\"\n//   3. Move the (cosmetic) \"_GoBack\" last-edit-position bookmark from the\n//      end of the authors paragraph to inside \"Albert Xue\" in the\n//      \"Author information\" section, mirroring where Word leaves it after\n//      the two text edits above.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// --- 1. \"public\" -> \"manuscript\" -------------------------------------\nconst publicResults = body.search(\"public\", { matchCase: true });\npublicResults.load(\"text\");\nawait context.sync();\n\nif (publicResults.items.length > 0) {\n  publicResults.items[0].insertText(\"manuscript\", \"Replace\");\n  await context.sync();\n}\n\n// --- 2. date \"02/06/17\" -> \"04/27/18\" ----------------------------------\nconst dateResults = body.search(\"02/06/17\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"04/27/18\", \"Replace\");\n  await context.sync();\n}\n\n// --- 3. relocate the \"_GoBack\" bookmark --------------------------------\nconst goBack = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBack.load(\"isNullObject\");\nawait context.sync();\n\nif (!goBack.isNullObject) {\n  doc.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// \"Albert Xue of the Bagheri Lab\" (Author information section) is unique,\n// unlike \"Albert Xue\" alone which also occurs earlier in the doc.\nconst authorInfoResults = body.search(\"Albert Xue of the Bagheri Lab\", { matchCase: true });\nauthorInfoResults.load(\"text\");\nawait context.sync();\n\nif (authorInfoResults.items.length > 0) {\n  const albertXResults = authorInfoResults.items[0].search(\"Albert X\", { matchCase: true });\n  albertXResults.load(\"text\");\n  await context.sync();\n\n  if (albertXResults.items.length > 0) {\n    const endRange = albertXResults.items[0].getRange(\"End\");\n    endRange.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Tiny update to readme:\n#   1. \"...alongside our public with the following authors:\"\n#        -> \"...alongside our manuscript with the following authors:\"\n#   2. \"Version 0.1 updated on 02/06/17\" -> \"Version 0.1 updated on 04/27/18\"\n#   3. Move the (cosmetic) \"_GoBack\" last-edit-position bookmark from the end\n#      of the authors paragraph to inside \"Albert Xue\" in the \"Author\n#      information\" section, mirroring where Word leaves it after the two\n#      text edits above.\n\n$d = $word.ActiveDocument\n\n# --- 1. \"public\" -> \"manuscript\" --------------------------------------\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.Text = \"public\"\nif ($find1.Execute()) {\n    $range1.Text = \"manuscript\"\n}\n\n# --- 2. date \"02/06/17\" -> \"04/27/18\" -----------------------------------\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.Text = \"02/06/17\"\nif ($find2.Execute()) {\n    $range2.Text = \"04/27/18\"\n}\n\n# --- 3. relocate the \"_GoBack\" bookmark ---------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.Text = \"Albert Xue of the Bagheri Lab\"\nif ($find3.Execute()) {\n    # \"Albert X\" is 8 characters; put a collapsed bookmark right after it.\n    $point = $d.Range($range3.Start + 8, $range3.Start + 8)\n    $d.Bookmarks.Add(\"_GoBack\", $point)\n}\n"}
